{"js": "// Turns the standalone paragraph \"neue Start\" into \"ein neuer Start beginnt \"\n// by typing at three caret positions, the same way the author produced the\n// change (the saved OOXML shows the inserted text split across separate\n// runs rather than merged into the pre-existing \"neue\" / \" Start\" runs).\n//\n// We re-search before each insertion so every step works off a fresh,\n// correctly-positioned range instead of a stale one.\n\nconst body = context.document.body;\n\n// 1) Put the caret right before \"neue Start\" and type \"ein \".\nlet results = body.search(\"neue Start\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].getRange(\"Start\").insertText(\"ein \", \"Before\");\n  await context.sync();\n\n  // 2) Put the caret right after \"neue\" (before \" Start\") and type \"r\",\n  //    turning \"neue\" into \"neuer\".\n  results = body.search(\"neue\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].getRange(\"End\").insertText(\"r\", \"Before\");\n  await context.sync();\n\n  // 3) Put the caret right after \"neuer Start\" and type \" beginnt \".\n  results = body.search(\"neuer Start\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].getRange(\"End\").insertText(\" beginnt \", \"Before\");\n  await context.sync();\n}\n", "ps1": "# Turns the standalone paragraph \"neue Start\" into \"ein neuer Start beginnt \"\n# by placing the caret at three positions and typing, mirroring how the\n# author produced the change (the saved OOXML shows the inserted text split\n# across separate runs rather than merged into the pre-existing \"neue\" /\n# \" Start\" runs).\n#\n# Each step re-runs Find against $d.Content so it always works off a fresh,\n# correctly-positioned range instead of a stale one.\n\n$d = $word.ActiveDocument\n\n# 1) Put the caret right before \"neue Start\" and type \"ein \".\n$rng1 = $d.Content\nif ($rng1.Find.Execute(\"neue Start\")) {\n    $insPoint1 = $d.Range($rng1.Start, $rng1.Start)\n    $insPoint1.InsertBefore(\"ein \")\n\n    # 2) Put the caret right after \"neue\" (before \" Start\") and type \"r\",\n    #    turning \"neue\" into \"neuer\".\n    $rng2 = $d.Content\n    $rng2.Find.Execute(\"neue\") | Out-Null\n    $insPoint2 = $d.Range($rng2.End, $rng2.End)\n    $insPoint2.InsertBefore(\"r\")\n\n    # 3) Put the caret right after \"neuer Start\" and type \" beginnt \".\n    $rng3 = $d.Content\n    $rng3.Find.Execute(\"neuer Start\") | Out-Null\n    $insPoint3 = $d.Range($rng3.End, $rng3.End)\n    $insPoint3.InsertBefore(\" beginnt \")\n}\n"}
